$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.670.12'
$ws.Range("E2").Value = '  +0.60%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.861.12'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '524.35'
$ws.Range("E5").Value = '  +7.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.26'
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  -2.30%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("E9").Value = '  -4.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  -5.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000325'
$ws.Range("E11").Value = '  -6.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '41.65'
$ws.Range("E12").Value = '  -2.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.464.43'
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("E14").Value = '  -3.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.88'
$ws.Range("E15").Value = '  +9.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.867.04'
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("E17").Value = '  +6.63%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.99'
$ws.Range("E18").Value = '  -1.10%  '
$ws.Range("E19").Value = '  -1.51%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.751.96'
$ws.Range("E20").Value = '  +0.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '419.57'
$ws.Range("E21").Value = '  -3.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.39'
$ws.Range("E22").Value = '  -3.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.99'
$ws.Range("E23").Value = '  -4.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.93'
$ws.Range("E24").Value = '  -4.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.97'
$ws.Range("E25").Value = '  +5.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.28'
$ws.Range("E26").Value = '  -9.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.58'
$ws.Range("E27").Value = '  -4.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '35.90'
$ws.Range("E28").Value = '  -3.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '682.30'
$ws.Range("E29").Value = '  -4.55%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.09'
$ws.Range("E30").Value = '  -1.92%  '
$ws.Range("E31").Value = '  -3.50%  '
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '67.13'
$ws.Range("E33").Value = '  +9.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.434'
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("E35").Value = '  -4.69%  '
$ws.Range("D36").Value = ('0.0{0}0844' -f [char]0x2083)
$ws.Range("E36").Value = '  -3.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '39.69'
$ws.Range("E37").Value = '  -2.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.147'
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.19'
$ws.Range("E41").Value = '  +1.91%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0477'
$ws.Range("E42").Value = '  -3.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.18'
$ws.Range("E43").Value = '  +4.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.71'
$ws.Range("E44").Value = '  -12.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.41'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.780.41'
$ws.Range("E47").Value = '  +14.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.93'
$ws.Range("E48").Value = '  +4.75%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = ('0.0{0}0344' -f [char]0x2086)
$ws.Range("E49").Value = '  -8.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '144.43'
$ws.Range("E50").Value = '  +1.31%  '
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000266'
$ws.Range("E51").Value = '  +12.01%  '
